$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.412
$ws.Range("B4").Value = 5.474000000000001
$ws.Range("A6").Value = -22.087
$ws.Range("A7").Value = -20.637
$ws.Range("A8").Value = -21.779
$ws.Range("B8").Value = 5.47
$ws.Range("B9").Value = 5.808
$ws.Range("B12").Value = 5.52
$ws.Range("A16").Value = -21.798
$ws.Range("B17").Value = 5.153
$ws.Range("B18").Value = 5.87
$ws.Range("B19").Value = 7.491
$ws.Range("A20").Value = -21.902
$ws.Range("B20").Value = 5.289
$ws.Range("A21").Value = -20.186
$ws.Range("B26").Value = 6.291
$ws.Range("A28").Value = -21.624
$ws.Range("A29").Value = -21.5
$ws.Range("A30").Value = -21.513
$ws.Range("B31").Value = 6.15
$ws.Range("A32").Value = -21.418
$ws.Range("B39").Value = 6.770999999999999
$ws.Range("A40").Value = -20.624
$ws.Range("B40").Value = 7.25
$ws.Range("B41").Value = 6.461999999999999
$ws.Range("B42").Value = 6.223999999999999
$ws.Range("B43").Value = 6.071000000000001
$ws.Range("A46").Value = -21.489
$ws.Range("B47").Value = 6.263999999999999
$ws.Range("B48").Value = 5.443
$ws.Range("A51").Value = -21.176
$ws.Range("A52").Value = -21.646
$ws.Range("B54").Value = 5.313000000000001
$ws.Range("A57").Value = -21.889
$ws.Range("A59").Value = -22.257
$ws.Range("A62").Value = -21.897
$ws.Range("B62").Value = 5.336
$ws.Range("B63").Value = 5.252
$ws.Range("B64").Value = 5.526999999999999
$ws.Range("A66").Value = -21.563
$ws.Range("A73").Value = -21.367
$ws.Range("A74").Value = -20.673
$ws.Range("B76").Value = 6.114
$ws.Range("A77").Value = -21.503
$ws.Range("B81").Value = 5.441
$ws.Range("B84").Value = 5.924000000000001
$ws.Range("B89").Value = 5.332000000000001
$ws.Range("A92").Value = -21.566
$ws.Range("B94").Value = 5.795
$ws.Range("A100").Value = -22.031
